$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3894767.96
$ws.Range("C9").Value = 616550.24
$ws.Range("D9").Value = 4511318.2
$ws.Range("E9").Value = 13.66674246121677
$ws.Range("F9").Value = 86.33325753878323
$ws.Range("G9").Value = -40.41349356522519
$ws.Range("H9").Value = -29.66587414496372
$ws.Range("I9").Value = 39580
$ws.Range("J9").Value = 1694
$ws.Range("K9").Value = 41274
$ws.Range("L9").Value = 28559
$ws.Range("M9").Value = 157.9648517104941
$ws.Range("N9").Value = 7.845821648110207
